$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: rename B1, keep C1:G1 (1..5), drop H1:L1 ---
$ws.Range("B1").Value = "Values"
$ws.Range("H1:L1").Clear()

# --- Row 2: becomes the "Final Value" data row ---
$ws.Range("B2").Value = "Final Value"
$ws.Range("C2").Value = 1691505.018026276
$ws.Range("D2").Value = 1692665.373855845
$ws.Range("E2").Value = 1691182.308938605
$ws.Range("F2").Value = 1692047.795448294
$ws.Range("G2").Value = 1692222.519107871
$ws.Range("H2:L2").Clear()

# --- Row 3 (new): "Annualized Return" ---
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Annualized Return"
$ws.Range("C3").Value = 0.1895782454560733
$ws.Range("D3").Value = 0.189847701033854
$ws.Range("E3").Value = 0.1895032845774076
$ws.Range("F3").Value = 0.1897043035800674
$ws.Range("G3").Value = 0.1897448767658645

# --- Row 4 (new): "Sharpe Ratio" ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Sharpe Ratio"
$ws.Range("C4").Value = 0.7406116551948297
$ws.Range("D4").Value = 0.7417213538127997
$ws.Range("E4").Value = 0.7402971903190578
$ws.Range("F4").Value = 0.7411251441401773
$ws.Range("G4").Value = 0.7410931340116486

$excel.CutCopyMode = $false
